$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JdT-TPI_LRD")

# Remove the leftover SUM formula in H7 (no longer needed)
$ws.Range("H7").ClearContents()

# Add the two new work-diary rows (52 and 53)
$ws.Range("A52").Value = 44694
$ws.Range("B52").Value = "Réalisation"
$ws.Range("C52").Value = 0.75
$ws.Range("D52").Value = "Revue des commentaires, définitions de fonctions, etc"

$ws.Range("A53").Value = 44694
$ws.Range("B53").Value = "Analyse"
$ws.Range("C53").Value = 0.75
$ws.Range("D53").Value = "Reprise de la documentation du projet"
$ws.Range("E53").Value = "Pas beaucoup de documentation cette semaine"

# Match the date number format used for the other recent rows (same style as A44:A51)
$ws.Range("A44").Copy()
$ws.Range("A52:A53").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A52").Value = 44694
$ws.Range("A53").Value = 44694

# Resize the table / list object so it covers the new rows
$table = $ws.ListObjects.Item("Tableau1")
$table.Resize($ws.Range("A1:F53"))

# Restore the view: select H7 (also drops the stale top-left scroll offset)
$ws.Range("H7").Select()
